$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "ldjf"
$ws.Range("C6").Value = "kdjkdjf"
$ws.Range("C7").Value = "dlff"
$ws.Range("F5").Value = "ijdf"

$ws.Range("G10").Select()
